# Edit: merge the "L'utilisateur est inscrit sur le site" precondition into
# the "L'utilisateur est connecté" precondition (as "... connecté sur le
# site"), and move the stray "_GoBack" bookmark (previously sitting in the
# middle of the "N11" sentence) to the end of that merged paragraph.

$d = $word.ActiveDocument

# --- Step 1: remove the old "_GoBack" bookmark (currently inside the N11
# sentence, between "1" and " : Le système envoie un mail...") -------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 2: locate the paragraph whose text is exactly
# "L'utilisateur est inscrit sur le site" ------------------------------------
$idxInscrit = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq "L’utilisateur est inscrit sur le site`r") {
        $idxInscrit = $i
    }
}

# --- Step 3: delete that whole paragraph (text + paragraph mark) -----------
$pInscrit = $d.Paragraphs.Item($idxInscrit)
$pInscrit.Range.Delete()

# --- Step 4: the following paragraph ("L'utilisateur est connecté") now
# sits at the same index; append " sur le site" to its text ----------------
$pConnecte = $d.Paragraphs.Item($idxInscrit)
$pConnecte.Range.InsertAfter(" sur le site")

# --- Step 5: add a fresh "_GoBack" bookmark collapsed at the very end of
# that paragraph's text (right before the paragraph mark). The COM host
# mis-places a bookmark collapsed exactly on "paragraph end - 1", so work
# around it: append a throw-away character, bookmark right before it, then
# delete the throw-away character. Freshly-built $d.Range(start,end) objects
# are used throughout (rather than .Duplicate()/.Collapse() chains) since
# those derived ranges are handled less reliably by the COM host. ----------
$insertPos = $pConnecte.Range.End - 1
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertAfter("X")

$pEnd = $pConnecte.Range.End
$markPos = $pEnd - 2
$markRange = $d.Range($markPos, $markPos)
$d.Bookmarks.Add("_GoBack", $markRange)

$tempCharRange = $d.Range($pEnd - 2, $pEnd - 1)
$tempCharRange.Delete()
